# SE 2018 test1 score commit
# Fills in Test 1 scores (column F) on the "SE" sheet for rows 4-51,
# corrects the submission date in F3 on both "SE" and "NE" sheets, and
# updates the active sheet/selection view state to match.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# SE sheet: enter Test 1 (column F) scores and the date header
# ---------------------------------------------------------------
$se = $wb.Worksheets.Item("SE")
$se.Activate()

# Date header (was stored as the shared-string "dd/mm/yy", now a real date value)
$se.Range("F3").Value = 44096

# Test 1 scores for rows 4-51, in row order (row 4 first)
$seScores = @(
    0, 8, 8, 10, 10, 10, 9, 9, 9, 10,
    10, 8, 9, 8, 9, 9, 10, 9, 10, 9,
    10, 9, 10, 10, 9, 10, 8, 9, 8, 9,
    10, 10, 9, 10, 9, 10, 8, 8, 9, 10,
    8, 10, 10, 10, 10, 7, 10, 8
)

$startRow = 4
for ($i = 0; $i -lt $seScores.Count; $i++) {
    $se.Cells.Item($startRow + $i, 6).Value = $seScores[$i]
}

$excel.ActiveWindow.Zoom = 80
$se.Range("F50").Select()

# ---------------------------------------------------------------
# NE sheet: fix the submission date back to 44096 and reset the view
# ---------------------------------------------------------------
$ne = $wb.Worksheets.Item("NE")
$ne.Activate()
$ne.Range("F3").Value = 44096
$ne.Range("F3").Select()

# ---------------------------------------------------------------
# Leave SE as the active/selected sheet
# ---------------------------------------------------------------
$se.Activate()
$se.Range("F50").Select()
